$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellAddr, $Text)
    $escaped = $Text -replace '"', '""'
    $ws.Range("ZZ1").Formula = '="' + $escaped + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($CellAddr).PasteSpecial(-4163)
    $ws.Range("ZZ1").ClearContents()
}

Set-TextValue "D2" '30.107.41'
$ws.Range("E2").Value = '  +0.22%  '

Set-TextValue "D3" '1.926.46'
$ws.Range("E3").Value = '  +2.85%  '

Set-TextValue "D4" '1.002'
$ws.Range("E4").Value = '  +0.06%  '

Set-TextValue "D5" '320.16'
$ws.Range("E5").Value = '  +0.19%  '

Set-TextValue "D6" '1.001'
$ws.Range("E6").Value = '  +0.05%  '

Set-TextValue "D7" '0.5080'
$ws.Range("E7").Value = '  +0.48%  '

Set-TextValue "D8" '0.4028'
$ws.Range("E8").Value = '  +2.18%  '

Set-TextValue "D9" '0.08352'
$ws.Range("E9").Value = '  +1.76%  '

Set-TextValue "D10" '1.118'
$ws.Range("E10").Value = '  +2.27%  '

Set-TextValue "D11" '42.20'
$ws.Range("E11").Value = '  -0.03%  '

Set-TextValue "D12" '24.04'
$ws.Range("E12").Value = '  +0.82%  '

Set-TextValue "D13" '6.413'
$ws.Range("E13").Value = '  +1.71%  '

Set-TextValue "D14" '1.917.97'
$ws.Range("E14").Value = '  +2.72%  '

Set-TextValue "D15" '7.277'
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("E16").Value = '  +0.01%  '

Set-TextValue "D17" '92.87'
$ws.Range("E17").Value = '  +0.80%  '

Set-TextValue "D18" '0.00001097'
$ws.Range("E18").Value = '  +0.55%  '

Set-TextValue "D19" '0.06509'
$ws.Range("E19").Value = '  +1.34%  '

$ws.Range("E20").Value = '  +2.17%  '

Set-TextValue "D21" '1.000'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("E22").Value = '  +2.13%  '

Set-TextValue "D23" '30.110.64'
$ws.Range("E23").Value = '  +0.26%  '

Set-TextValue "D24" '11.38'
$ws.Range("E24").Value = '  +2.14%  '

Set-TextValue "D25" '2.187'
$ws.Range("E25").Value = '  +0.70%  '

Set-TextValue "D26" '2.138.95'
$ws.Range("E26").Value = '  +2.56%  '

Set-TextValue "D27" '21.91'
$ws.Range("E27").Value = '  +2.37%  '

Set-TextValue "D28" '162.64'
$ws.Range("E28").Value = '  +1.52%  '

Set-TextValue "D29" '2.273'
$ws.Range("E29").Value = '  +1.75%  '

Set-TextValue "D30" '129.06'
$ws.Range("E30").Value = '  +1.21%  '

Set-TextValue "D31" '1.140'
$ws.Range("E31").Value = '  +7.02%  '

$ws.Range("E32").Value = '  +1.18%  '

Set-TextValue "D33" '5.970'
$ws.Range("E33").Value = '  +0.26%  '

Set-TextValue "D34" '3.791'
$ws.Range("E34").Value = '  +2.82%  '

Set-TextValue "D35" '0.02455'
$ws.Range("E35").Value = '  +1.01%  '

Set-TextValue "D36" '5.323'
$ws.Range("E36").Value = '  +1.56%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D37" '1.258'
$ws.Range("E37").Value = '  +6.85%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D38" '0.06466'
$ws.Range("E38").Value = '  +1.32%  '

Set-TextValue "D39" '0.2153'
$ws.Range("E39").Value = '  +0.19%  '

Set-TextValue "D40" '0.6488'
$ws.Range("E40").Value = '  +2.70%  '

Set-TextValue "D41" '8.683'
$ws.Range("E41").Value = '  +2.03%  '

Set-TextValue "D42" '11.67'
$ws.Range("E42").Value = '  +2.46%  '

Set-TextValue "D43" '1.216'
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D44" '0.6066'
$ws.Range("E44").Value = '  +2.37%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D45" '13.29'
$ws.Range("E45").Value = '  +2.83%  '

Set-TextValue "D46" '2.167'
$ws.Range("E46").Value = '  +4.15%  '

Set-TextValue "D47" '3.625'
$ws.Range("E47").Value = '  -0.23%  '

Set-TextValue "D48" '122.43'
$ws.Range("E48").Value = '  -0.56%  '

$ws.Range("E49").Value = '  +0.13%  '

Set-TextValue "D50" '1.131'
$ws.Range("E50").Value = '  +0.93%  '

Set-TextValue "D51" '77.97'
$ws.Range("E51").Value = '  +1.01%  '
